# Update the underlying score inputs on the sheet.
# The table in columns L and S holds the raw "score" text values that feed
# formulas elsewhere (M:Q and T:X), and B5 holds a team abbreviation used
# by formulas in column E/H. Updating these source cells lets Excel
# recalculate all of the dependent formulas/data tables automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (team abbreviation input)
$ws.Range("B5").Value = "LAR"

# Column L (left-side "W ex R" score strings)
$ws.Range("L2").Value  = "  27   65"
$ws.Range("L3").Value  = "  20   57"
$ws.Range("L4").Value  = "  24   52"
$ws.Range("L5").Value  = "  31   50"
$ws.Range("L6").Value  = "  30   38"
$ws.Range("L7").Value  = "  23   36"
$ws.Range("L8").Value  = "  34   35"
$ws.Range("L9").Value  = "  17   31"
$ws.Range("L10").Value = "  28   28"
$ws.Range("L11").Value = "  37   27"
$ws.Range("L12").Value = "  19   26"

# Column S (right-side "L ex R" score strings)
$ws.Range("S2").Value  = "  17   79"
$ws.Range("S3").Value  = "  10   74"
$ws.Range("S4").Value  = "  16   54"
$ws.Range("S5").Value  = "  13   50"
$ws.Range("S6").Value  = "  20   47"
$ws.Range("S7").Value  = "  14   40"
$ws.Range("S8").Value  = "  24   36"
$ws.Range("S9").Value  = "  3    31"
$ws.Range("S10").Value = "  21   27"
$ws.Range("S11").Value = "  9    26"
$ws.Range("S12").Value = "  6    24"

$wb.Save()
